# Update Inventory Master for Feb NMPS Movies - Added February Movies
#
# This reproduces the author's edit: sixteen new movie rows were appended
# to the bottom of the "Sheet1" table (in the order the titles were typed),
# then that newly-typed block was sorted alphabetically by TITLE and the
# running ID numbers (column A) were refreshed to stay sequential with the
# row position - matching the existing convention used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$app = $ws.Application
$win = $app.ActiveWindow
$lo  = $ws.ListObjects.Item(1)

# Titles + expiration dates, in the order they were originally entered
# (this is also the order the new strings land in the shared-strings table).
$entries = @(
    @{Title="THE GOOD DINOSAUR";                       Date=43800},
    @{Title="STAR WARS: THE FORCE AWAKENS";             Date=43800},
    @{Title="POINT BREAK";                              Date=43831},
    @{Title="JEM AND THE HOLOGRAMS";                    Date=43770},
    @{Title="CRIMSON PEAK";                             Date=43739},
    @{Title="VICTOR FRANKENSTEIN";                      Date=43800},
    @{Title="TRUTH";                                    Date=43800},
    @{Title="STEVE JOBS";                                Date=43770},
    @{Title="MISS YOU ALREADY";                         Date=43800},
    @{Title="DADDY'S HOME";                             Date=43800},
    @{Title="SPOTLIGHT";                                Date=43800},
    @{Title="AIRPLANE";                                 Date=43344},
    @{Title="IN THE HEART OF THE SEA";                  Date=43800},
    @{Title="BROOKLYN";                                 Date=43800},
    @{Title="ALVIN AND THE CHIPMUNKS: THE ROAD CHIP";   Date=43800},
    @{Title="KRAMPUS";                                  Date=43800}
)

$startRow = 712
$count    = $entries.Count
$endRow   = $startRow + $count - 1

# Write the new rows at the bottom of the sheet, in entry order, with
# sequential IDs just like every other append in this workbook.
for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 711 + $i
    $ws.Cells.Item($r, 2).Value = $entries[$i].Title
    $ws.Cells.Item($r, 3).Value = $entries[$i].Date
    $ws.Cells.Item($r, 3).NumberFormat = "mmm-yy"
}

# Sort just the newly added block alphabetically by TITLE, as was done
# for every previous batch of additions in this sheet.
$newBlock = $ws.Range("A" + $startRow + ":C" + $endRow)
$sortKey  = $ws.Range("B" + $startRow + ":B" + $endRow)
$newBlock.Sort($sortKey, 1)

# Re-sequence the ID column so it again lines up with row position.
for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 711 + $i
}

# Grow the table (Table1) so it covers the new rows.
$lo.Resize($ws.Range("A1:C" + $endRow))

# Re-apply the freeze (top row) / scroll position / selection the author
# left the sheet in after adding the new titles.
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.Zoom = 100
$ws.Range("A692").Select()
$ws.Rows.Item($startRow).Select()
